# Update "Handback" report timestamps to reflect the new generation run.
$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G)
$overview.Range("G2").Value = "2016-08-20 01:07:23"

# zh-cn sheet: "Correspond Handoff Datetime" (column H) and
# "Correspond Handback DateTime" (column K)
$zhcn.Range("H2").Value = "2016-08-20 01:07:18"
$zhcn.Range("K2").Value = "2016-08-20 01:07:35"

# de-de sheet: "Correspond Handoff Datetime" (column H) and
# "Correspond Handback DateTime" (column K)
$dede.Range("H2").Value = "2016-08-20 01:07:23"
$dede.Range("K2").Value = "2016-08-20 01:07:41"
